$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.209705352783203
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = 3.929324626922607
$ws.Range("D1").Value = 2.01487922668457
$ws.Range("E1").Value = 1.444768190383911
